$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "AddVisitReportAndScheduleForSite"
$ws.Range("A5").Select()
